# The slide 2 body placeholder asks students to insert a visualization of
# the Gaussian kernel. Update the question so it explicitly asks for both a
# 1D and a 2D visualization.
$p = $ppt.ActivePresentation

$shape = $null
foreach ($sl in $p.Slides) {
    foreach ($sh in $sl.Shapes) {
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -like "*Gaussian kernel*") {
                $shape = $sh
            }
        }
    }
}

$tr = $shape.TextFrame.TextRange

# Re-assert the original prompt text (content unchanged) and touch the space
# in the middle of the sentence so it becomes its own run, matching how
# PowerPoint splits a run when text already on the slide is nudged/re-edited.
$tr.Text = "[insert visualization of Gaussian kernel from project-1.ipynb here]"
$mid = $tr.Characters(41, 1)
$mid.Text = " "

# Append the new "1D:" / "2D:" sub-question paragraphs, with a blank line
# separating them from the prompt and from each other.
$null = $tr.InsertAfter("`r1D: `r`r2D:")
